# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the "Weekly Quantity" sheet
#    to "Weekly_PO_Qty".
# 2. Rename the "Requested quantity" header on the "Monthly Trend" sheet
#    to "Monthly_PO_Qty".
# 3. Add a new "PO Forecast" worksheet (placed after "Monthly Trend") with
#    forecast data: ds, PO_Forecast, yhat_lower, yhat_upper.

$wb = $excel.ActiveWorkbook

# --- Step 1: Weekly Quantity sheet header rename ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Step 2: Monthly Trend sheet header rename ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 3: Add PO Forecast sheet (placed after "Monthly Trend", the
#     current last tab) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used by the rest of the workbook (0.75/0.75/1/1/.5/.5 in)
$wsForecast.PageSetup.LeftMargin = $wsWeekly.PageSetup.LeftMargin
$wsForecast.PageSetup.RightMargin = $wsWeekly.PageSetup.RightMargin
$wsForecast.PageSetup.TopMargin = $wsWeekly.PageSetup.TopMargin
$wsForecast.PageSetup.BottomMargin = $wsWeekly.PageSetup.BottomMargin
$wsForecast.PageSetup.HeaderMargin = $wsWeekly.PageSetup.HeaderMargin
$wsForecast.PageSetup.FooterMargin = $wsWeekly.PageSetup.FooterMargin

# Seed the header row + date column by copying the existing bold/bordered
# header style and the date-formatted style from "Weekly Quantity" so the
# new sheet reuses the workbook's existing cell styles instead of minting
# new ones.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("C1:D1"))
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A17"))

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(44934.99999999999, 4, 0.2583725556269224, 8.544720530014688),
    @(44941.99999999999, 5, 0.4744135838433373, 8.459683813688761),
    @(44948.99999999999, 5, 0.838795761777651,  8.854403084947119),
    @(44955.99999999999, 5, 1.073755002169896,  8.643154520314495),
    @(44969.99999999999, 5, 1.161409841861742,  9.055563190976011),
    @(44976.99999999999, 5, 1.374381303212875,  9.053337001813816),
    @(44983.99999999999, 5, 1.409785683330566,  9.338109641644007),
    @(44990.99999999999, 6, 1.627990593878511,  9.509461279270129),
    @(44997.99999999999, 6, 1.767326365728237,  9.838450369148669),
    @(45004.99999999999, 6, 2.044153586669379,  9.788107433769172),
    @(45011.99999999999, 6, 2.132750910706613,  9.596873652986835),
    @(45018.99999999999, 6, 2.11119452452593,   10.18855843018231),
    @(45025.99999999999, 6, 2.187906623151452,  10.13073544007524),
    @(45032.99999999999, 6, 2.448532865279946,  9.877744429978328),
    @(45039.99999999999, 7, 2.783550391812425,  10.2634823219367),
    @(45046.99999999999, 7, 2.786878409781183,  10.61460160003905)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
